$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2265625
$ws.Range("C2").Value = 0.4765625
$ws.Range("J2").Value = 0.01953125
$ws.Range("O2").Value = 0.00390625
$ws.Range("P2").Value = 0.15625
$ws.Range("S2").Value = 0.1171875

# Row 3
$ws.Range("B3").Value = 0.01626016260162602
$ws.Range("C3").Value = 0.04065040650406504
$ws.Range("J3").Value = 0.008130081300813009
$ws.Range("P3").Value = 0.6991869918699187
$ws.Range("S3").Value = 0.2357723577235772

# Row 4
$ws.Range("J4").Value = 0.05882352941176471
$ws.Range("P4").Value = 0.6764705882352942
$ws.Range("S4").Value = 0.2647058823529412

# Row 6
$ws.Range("B6").Value = 0.0625
$ws.Range("D6").Value = 0.02604166666666667
$ws.Range("E6").Value = 0.005208333333333333
$ws.Range("F6").Value = 0.05208333333333334
$ws.Range("J6").Value = 0.2864583333333333
$ws.Range("O6").Value = 0.02083333333333333
$ws.Range("Q6").Value = 0.1197916666666667
$ws.Range("R6").Value = 0.07291666666666667
$ws.Range("S6").Value = 0.3541666666666667

# Row 7
$ws.Range("B7").Value = 0.08490566037735849
$ws.Range("D7").Value = 0.01415094339622642
$ws.Range("F7").Value = 0.0660377358490566
$ws.Range("J7").Value = 0.09905660377358491
$ws.Range("O7").Value = 0.02830188679245283
$ws.Range("Q7").Value = 0.1981132075471698
$ws.Range("R7").Value = 0.07547169811320754
$ws.Range("S7").Value = 0.4339622641509434

# Row 8
$ws.Range("B8").Value = 0.04749340369393139
$ws.Range("D8").Value = 0.01846965699208443
$ws.Range("E8").Value = 0.002638522427440633
$ws.Range("F8").Value = 0.0633245382585752
$ws.Range("J8").Value = 0.09762532981530343
$ws.Range("O8").Value = 0.02110817941952507
$ws.Range("Q8").Value = 0.1952506596306069
$ws.Range("R8").Value = 0.09234828496042216
$ws.Range("S8").Value = 0.4617414248021108

# Row 9
$ws.Range("B9").Value = 0.1116071428571429
$ws.Range("D9").Value = 0.008928571428571428
$ws.Range("E9").Value = 0.004464285714285714
$ws.Range("F9").Value = 0.05357142857142857
$ws.Range("J9").Value = 0.1071428571428571
$ws.Range("O9").Value = 0.01339285714285714
$ws.Range("Q9").Value = 0.1875
$ws.Range("R9").Value = 0.05803571428571429
$ws.Range("S9").Value = 0.4553571428571428

# Row 10
$ws.Range("B10").Value = 0.0982940698619009
$ws.Range("D10").Value = 0.01462225832656377
$ws.Range("E10").Value = 0.0008123476848090983
$ws.Range("F10").Value = 0.05605199025182778
$ws.Range("J10").Value = 0.1186027619821284
$ws.Range("O10").Value = 0.01299756295694557
$ws.Range("Q10").Value = 0.2266450040617384
$ws.Range("R10").Value = 0.07554833468724614
$ws.Range("S10").Value = 0.3964256701868399

# Row 11
$ws.Range("G11").Value = 0.1594202898550725
$ws.Range("J11").Value = 0.08405797101449275
$ws.Range("K11").Value = 0.2144927536231884
$ws.Range("L11").Value = 0.5333333333333333
$ws.Range("S11").Value = 0.008695652173913044

# Row 12
$ws.Range("G12").Value = 0.6963350785340314
$ws.Range("J12").Value = 0.2094240837696335
$ws.Range("K12").Value = 0.01047120418848168
$ws.Range("L12").Value = 0.03141361256544502
$ws.Range("S12").Value = 0.05235602094240838

# Row 13
$ws.Range("G13").Value = 0.7272727272727273
$ws.Range("J13").Value = 0.2045454545454546
$ws.Range("S13").Value = 0.06818181818181818

# Row 14
$ws.Range("G14").Value = 0.75
$ws.Range("J14").Value = 0.25

# Row 15
$ws.Range("F15").Value = 0.01648351648351648
$ws.Range("H15").Value = 0.1153846153846154
$ws.Range("I15").Value = 0.08241758241758242
$ws.Range("J15").Value = 0.3461538461538461
$ws.Range("K15").Value = 0.1043956043956044
$ws.Range("M15").Value = 0.005494505494505495
$ws.Range("N15").Value = 0.005494505494505495
$ws.Range("O15").Value = 0.01098901098901099
$ws.Range("S15").Value = 0.3131868131868132

# Row 16
$ws.Range("F16").Value = 0.01342281879194631
$ws.Range("H16").Value = 0.2281879194630873
$ws.Range("I16").Value = 0.1342281879194631
$ws.Range("J16").Value = 0.3691275167785235
$ws.Range("K16").Value = 0.1140939597315436
$ws.Range("M16").Value = 0.006711409395973154
$ws.Range("O16").Value = 0.006711409395973154
$ws.Range("S16").Value = 0.1275167785234899

# Row 17
$ws.Range("F17").Value = 0.03036876355748373
$ws.Range("H17").Value = 0.1540130151843818
$ws.Range("I17").Value = 0.07592190889370933
$ws.Range("J17").Value = 0.4338394793926247
$ws.Range("K17").Value = 0.1127982646420824
$ws.Range("M17").Value = 0.008676789587852495
$ws.Range("N17").Value = 0.004338394793926247
$ws.Range("O17").Value = 0.05422993492407809
$ws.Range("S17").Value = 0.1258134490238612

# Row 18
$ws.Range("F18").Value = 0.01169590643274854
$ws.Range("H18").Value = 0.1812865497076023
$ws.Range("I18").Value = 0.1228070175438596
$ws.Range("J18").Value = 0.3684210526315789
$ws.Range("K18").Value = 0.08187134502923976
$ws.Range("M18").Value = 0.01754385964912281
$ws.Range("O18").Value = 0.03508771929824561
$ws.Range("S18").Value = 0.1812865497076023

# Row 19
$ws.Range("F19").Value = 0.01751713632901752
$ws.Range("H19").Value = 0.1706016755521706
$ws.Range("I19").Value = 0.1020563594821021
$ws.Range("J19").Value = 0.3777608530083778
$ws.Range("K19").Value = 0.1264280274181264
$ws.Range("M19").Value = 0.02665651180502666
$ws.Range("N19").Value = 0.001523229246001523
$ws.Range("O19").Value = 0.06321401370906321
$ws.Range("S19").Value = 0.1142421934501142
